# Move the "BYAC290" marker from B39 (Travis Winston) to B42 (Charles Dowswell),
# and update the sheet view scroll position / selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old cell and write the value into its new location.
$ws.Range("B39").ClearContents()
$ws.Range("B42").Value = "BYAC290"

# Scroll the view down and move the active selection to match the new edit location.
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("B39").Select()
